$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style from an existing header cell (e.g., E1) to the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

$data = @(
  @(0,0,0),
  @(0,0,0),
  @(1,0,0),
  @(1,0,0),
  @(0,0,0),
  @(0,0,0),
  @(1,0,0),
  @(0,0,0),
  @(0,0,0),
  @(0,0,0),
  @(0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = [bool]$data[$i][0]
  $ws.Cells.Item($row, 7).Value = [bool]$data[$i][1]
  $ws.Cells.Item($row, 8).Value = [bool]$data[$i][2]
}
